$wb = $excel.ActiveWorkbook

# --- Rename header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value2 = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value2 = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Cells.Item(1,1).Value2 = "ds"
$wsForecast.Cells.Item(1,2).Value2 = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value2 = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value2 = "yhat_upper"

# Copy header formatting (bold, border, centered) from an existing sheet's header
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- Data rows ---
$rows = @(
  "45333.99999999999,0,-402.8454862954759,309.6701108884287",
  "45340.99999999999,0,-373.3011455146004,318.7880274765889",
  "45347.99999999999,0,-364.9210802164425,338.2469988829304",
  "45354.99999999999,0,-343.3013451914923,361.7607160640896",
  "45361.99999999999,16,-340.9249397511163,351.6381059345043",
  "45368.99999999999,32,-315.9645350577208,365.6719631920936",
  "45375.99999999999,47,-302.2955082937316,381.7041346888378",
  "45382.99999999999,63,-277.5196071475429,416.2262534867315",
  "45389.99999999999,79,-266.722339940078,421.1442713064323",
  "45396.99999999999,95,-253.1133360662537,426.8394496723719",
  "45403.99999999999,111,-225.1633569260262,447.4315316069661",
  "45410.99999999999,127,-194.6476871421604,489.4412736989978",
  "45431.99999999999,175,-163.0013105001692,514.5478755721612",
  "45445.99999999999,206,-145.0650063769097,550.0731243912966",
  "45459.99999999999,238,-97.71888875694786,589.5010129753268",
  "45466.99999999999,254,-95.28266512597784,610.3442984463859",
  "45473.99999999999,270,-68.66814149650791,631.691301758442",
  "45480.99999999999,286,-81.08421964348962,633.0326173243287",
  "45487.99999999999,302,-10.88388892190486,682.972357810742",
  "45515.99999999999,365,14.11694095695514,730.6559758177079",
  "45522.99999999999,381,34.60975304281894,719.1193624909382",
  "45536.99999999999,413,67.29384192762014,749.1763294990433",
  "45550.99999999999,445,108.2142966896875,808.4713473146738",
  "45557.99999999999,461,103.8298791922681,821.0205420271506",
  "45564.99999999999,477,152.0070146274653,793.9350035251108",
  "45599.99999999999,556,216.0663037283911,890.832049313333",
  "45606.99999999999,572,231.29657114394,929.4580129256225",
  "45613.99999999999,588,245.3740649364219,953.3716649020749",
  "45620.99999999999,604,236.3578108148911,933.2377566321136",
  "45627.99999999999,620,261.0928479747151,969.0520208923049",
  "45634.99999999999,636,315.9228063489807,1002.689714105686",
  "45641.99999999999,652,302.8523908418792,1003.301572550641",
  "45648.99999999999,668,324.3884835032702,1028.98495638875",
  "45655.99999999999,683,339.9996125115034,1043.709251624747"
)

$r = 2
foreach ($line in $rows) {
  $parts = $line.Split(",")
  $wsForecast.Cells.Item($r,1).Value2 = [double]$parts[0]
  $wsForecast.Cells.Item($r,2).Value2 = [double]$parts[1]
  $wsForecast.Cells.Item($r,3).Value2 = [double]$parts[2]
  $wsForecast.Cells.Item($r,4).Value2 = [double]$parts[3]
  $r = $r + 1
}

# Copy the date-format style from an existing date cell onto the new ds column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A35").PasteSpecial(-4122)

$wsForecast.Range("A1").Select() | Out-Null
